$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9192272424697876
$ws.Range("B1").Value = 1.704797029495239
$ws.Range("D1").Value = 1.852428674697876
$ws.Range("E1").Value = 1.096906542778015
